$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: b.md row (row 3) moves from "Handed back" to "Ready for
# handoff" with a new handoff timestamp.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-28-19 08:28:12"

# ---------------------------------------------------------------------------
# zh-cn sheet: b.md row (row 3) gets a new handoff file/status/datetime, and
# its "Latest Handoff File" hyperlink now points at (and displays) the new
# handoff file name.
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-19 08:28:09"

# The COM hyperlink object model here only supports whole-collection
# mutation cleanly (per-item property writes/deletes leave stray duplicate
# nodes behind), so rebuild the sheet's hyperlinks collection in place with
# the updated display text for D3.
$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/97e1b4f2cb1b4b83b16fdcf798046363bf183853/e2e/a.md", "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/97e1b4f2cb1b4b83b16fdcf798046363bf183853/e2e/a.md", "", "", ".md")
$zhcn.Hyperlinks.Add($zhcn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c99eb845f0a07f3ceb556803006d18666cdcd04f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/32632a4d59bd500cc452d23fed61dcd9f29e72c4/e2e/a.md", "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b45ea2591b29f7651b26c70e09a7011b4cd8584b/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/97e1b4f2cb1b4b83b16fdcf798046363bf183853/e2e/b.md", "", "", "b.md")
$zhcn.Hyperlinks.Add($zhcn.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/97e1b4f2cb1b4b83b16fdcf798046363bf183853/e2e/b.md", "", "", ".md")
$zhcn.Hyperlinks.Add($zhcn.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c99eb845f0a07f3ceb556803006d18666cdcd04f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/32632a4d59bd500cc452d23fed61dcd9f29e72c4/e2e/a.md", "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b45ea2591b29f7651b26c70e09a7011b4cd8584b/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")

# ---------------------------------------------------------------------------
# de-de sheet: same shape of change as zh-cn, with de-de file names/links.
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-19 08:28:12"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/97e1b4f2cb1b4b83b16fdcf798046363bf183853/e2e/a.md", "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/97e1b4f2cb1b4b83b16fdcf798046363bf183853/e2e/a.md", "", "", ".md")
$dede.Hyperlinks.Add($dede.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/af93be5a9ad7cfc272a310c0f8e7d15b802e5fed/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/bc12dc84cf018f0f68e66385be6b05cd60c653fc/e2e/a.md", "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/99e1c0db84e5551b77d910d45fded78e834f2fcb/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/97e1b4f2cb1b4b83b16fdcf798046363bf183853/e2e/b.md", "", "", "b.md")
$dede.Hyperlinks.Add($dede.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/97e1b4f2cb1b4b83b16fdcf798046363bf183853/e2e/b.md", "", "", ".md")
$dede.Hyperlinks.Add($dede.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/af93be5a9ad7cfc272a310c0f8e7d15b802e5fed/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/bc12dc84cf018f0f68e66385be6b05cd60c653fc/e2e/a.md", "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/99e1c0db84e5551b77d910d45fded78e834f2fcb/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")

Write-Output "Report generated for handoff"
